$wb = $excel.ActiveWorkbook

# --- Step 1: add the new "2022-Q1" sheet -----------------------------------
# It shares the exact same layout/formatting as the "2021-Q4" sheet, so the
# simplest and most faithful way to create it is to copy that sheet and drop
# it right before "总计" (matching the new tab order: 2021-Q3, 2021-Q4,
# 2022-Q1, 总计).
$q4 = $wb.Worksheets.Item("2021-Q4")
$zongji_target = $wb.Worksheets.Item("总计")
$q4.Copy($zongji_target)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Sheet handles captured before a sheet-insertion resolve by tab position, so
# re-fetch "总计" by name now that it has shifted one slot to the right.
$zongji = $wb.Worksheets.Item("总计")

# Update the figures that differ from the 2021-Q4 snapshot. D:G are stored as
# text in the source workbook, so force text entry with a leading apostrophe
# and strip the resulting "quote prefix" style back to Normal so the cell
# matches the plain (unstyled) text cells used throughout the workbook.
$q1.Range("D2").Value = "'0.29"
$q1.Range("D2").Style = "Normal"
$q1.Range("E2").Value = "'94.14"
$q1.Range("E2").Style = "Normal"
$q1.Range("F2").Value = "'2.46"
$q1.Range("F2").Style = "Normal"
$q1.Range("G2").Value = "'0.0071"
$q1.Range("G2").Style = "Normal"
$q1.Range("H2").Value = 5

# --- Step 2: update the "总计" (summary) sheet ------------------------------
# A new top row for 2022-Q1 is inserted and the existing two rows shift down
# one position, with their running index (column A) bumped accordingly.
# Row 4 is brand new, so first clone row 3's formatting (the numbered-row
# style in column A) down onto it via a cell copy, then overwrite the values
# on every row.
$zongji.Range("A3").Copy($zongji.Range("A4"))

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2021-Q3"
$zongji.Range("C4").Value = 1
$zongji.Range("D4").Value = 0.01

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2021-Q4"
$zongji.Range("C3").Value = 1
$zongji.Range("D3").Value = 0.01

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 1
$zongji.Range("D2").Value = 0.01

# Restore the originally-active tab; the various Copy/rename steps above
# leave the newly touched sheet selected, but the source workbook had
# "2021-Q3" active and the diff doesn't touch that state.
$wb.Worksheets.Item("2021-Q3").Activate()
